$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 64486.57
$ws.Range("J62").Value = 64486.57
$ws.Range("L62").Value = 64486.57
$ws.Range("N62").Value = -65734.57000000001

$ws.Range("H65").Value = 64486.57
$ws.Range("J65").Value = 64486.57
$ws.Range("L65").Value = 322432.85
$ws.Range("N65").Value = -328672.85

$ws.Range("H74").Value = 150010980
$ws.Range("J74").Value = 19800.8
$ws.Range("L74").Value = 19800.8
$ws.Range("N74").Value = -21672.8

$ws.Range("H77").Value = 150010980
$ws.Range("J77").Value = 19800.8
$ws.Range("L77").Value = 99004
$ws.Range("N77").Value = -108364

$ws.Range("H80").Value = 32439.625
$ws.Range("I80").Value = 11469.889
$ws.Range("K80").Value = 34409.667
$ws.Range("M80").Value = -33411.667

$ws.Range("H83").Value = 32439.625
$ws.Range("I83").Value = 11469.889
$ws.Range("K83").Value = 103229.001
$ws.Range("M83").Value = -98237.00099999999

$ws.Range("H88").Value = 15185212
$ws.Range("I88").Value = 47622420
$ws.Range("J88").Value = 47849.2
$ws.Range("K88").Value = 47622420
$ws.Range("L88").Value = 47849.2
$ws.Range("M88").Value = -47622014
$ws.Range("N88").Value = -48661.2

$ws.Range("H91").Value = 15185212
$ws.Range("I91").Value = 47622420
$ws.Range("J91").Value = 47849.2
$ws.Range("K91").Value = 47622420
$ws.Range("L91").Value = 47849.2
$ws.Range("M91").Value = -47621016
$ws.Range("N91").Value = -50657.2

$ws.Range("H98").Value = 32260942
$ws.Range("I98").Value = 37039748
$ws.Range("K98").Value = 37039748
$ws.Range("M98").Value = -37038250

$ws.Range("H112").Value = 3482.568
$ws.Range("J112").Value = 3482.568
$ws.Range("L112").Value = 10447.704
$ws.Range("N112").Value = -12663.704

$ws.Range("H122").Value = 32260942
$ws.Range("I122").Value = 37039748
$ws.Range("K122").Value = 111119244
$ws.Range("M122").Value = -111116794

$ws.Range("H132").Value = 1294.0728
$ws.Range("I132").Value = 1177.804
$ws.Range("K132").Value = 3533.412
$ws.Range("M132").Value = -1003.412

$ws.Range("H137").Value = 2218.1614
$ws.Range("I137").Value = 1959.5306
$ws.Range("K137").Value = 5878.5918
$ws.Range("M137").Value = -3328.5918

$ws.Range("H138").Value = 1567830.9
$ws.Range("J138").Value = 2508055.8
$ws.Range("L138").Value = 7524167.399999999
$ws.Range("N138").Value = -7534447.399999999


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5389.97
$ws.Range("I32").Value = 4577.1445
$ws.Range("J32").Value = 31671.334
$ws.Range("K32").Value = 4577.1445
$ws.Range("L32").Value = 31671.334
$ws.Range("M32").Value = -4290.1445
$ws.Range("N32").Value = -32245.334

$ws.Range("H57").Value = 4974.25
$ws.Range("I57").Value = 4974.25
$ws.Range("K57").Value = 4974.25
$ws.Range("M57").Value = -4490.25

$ws.Range("H61").Value = 8426.200000000001
$ws.Range("I61").Value = 892.1111
$ws.Range("K61").Value = 892.1111
$ws.Range("M61").Value = -680.1111

$ws.Range("H122").Value = 3465.36
$ws.Range("I122").Value = 2365
$ws.Range("K122").Value = 7095
$ws.Range("M122").Value = -4645

$ws.Range("H136").Value = 8426.200000000001
$ws.Range("I136").Value = 892.1111
$ws.Range("K136").Value = 2676.3333
$ws.Range("M136").Value = -126.3332999999998

$ws.Range("H139").Value = 60592.145
$ws.Range("J139").Value = 60592.145
$ws.Range("L139").Value = 60592.145
$ws.Range("N139").Value = -70872.14499999999


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 43000
$ws.Range("J63").Value = 43000
$ws.Range("L63").Value = 43000
$ws.Range("N63").Value = -44372

$ws.Range("H66").Value = 43000
$ws.Range("J66").Value = 43000
$ws.Range("L66").Value = 129000
$ws.Range("N66").Value = -135864

$ws.Range("H113").Value = 5217
$ws.Range("I113").Value = 5217
$ws.Range("K113").Value = 5217
$ws.Range("M113").Value = -3047

$ws.Range("H134").Value = 4451.103
$ws.Range("I134").Value = 1639.2325
$ws.Range("J134").Value = 9287.52
$ws.Range("K134").Value = 4917.6975
$ws.Range("L134").Value = 27862.56
$ws.Range("M134").Value = -2382.6975
$ws.Range("N134").Value = -32932.56


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5058.253
$ws.Range("I31").Value = 2261.8215
$ws.Range("K31").Value = 2261.8215
$ws.Range("M31").Value = -1966.8215

$ws.Range("H34").Value = 5058.253
$ws.Range("I34").Value = 2261.8215
$ws.Range("K34").Value = 2261.8215
$ws.Range("M34").Value = -2059.8215

$ws.Range("H86").Value = 78130000
$ws.Range("I86").Value = 156250000
$ws.Range("K86").Value = 156250000
$ws.Range("M86").Value = -156248877

$ws.Range("H89").Value = 78130000
$ws.Range("I89").Value = 156250000
$ws.Range("K89").Value = 781250000
$ws.Range("M89").Value = -781244384

$ws.Range("H132").Value = 4790.6724
$ws.Range("I132").Value = 2050.4866
$ws.Range("J132").Value = 9618.619000000001
$ws.Range("K132").Value = 6151.459800000001
$ws.Range("L132").Value = 28855.857
$ws.Range("M132").Value = -3621.459800000001
$ws.Range("N132").Value = -33915.857

$ws.Range("H134").Value = 3416.551
$ws.Range("I134").Value = 1380.9706
$ws.Range("K134").Value = 4142.9118
$ws.Range("M134").Value = -1607.9118


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 3250
$ws.Range("J130").Value = 4000
$ws.Range("L130").Value = 12000
$ws.Range("N130").Value = -22040


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 59285.57
$ws.Range("I52").Value = 45000
$ws.Range("K52").Value = 45000
$ws.Range("M52").Value = -44741

$ws.Range("H80").Value = 2167.6667
$ws.Range("I80").Value = 2426.5
$ws.Range("K80").Value = 2426.5
$ws.Range("M80").Value = -1428.5

$ws.Range("H83").Value = 2167.6667
$ws.Range("I83").Value = 2426.5
$ws.Range("K83").Value = 12132.5
$ws.Range("M83").Value = -7140.5

$ws.Range("H102").Value = 5833.1665
$ws.Range("I102").Value = 5833.1665
$ws.Range("K102").Value = 5833.1665
$ws.Range("M102").Value = -4211.1665

$ws.Range("H132").Value = 5316.4414
$ws.Range("I132").Value = 2811.5107
$ws.Range("K132").Value = 8434.5321
$ws.Range("M132").Value = -5904.5321


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H53").Value = 9500
$ws.Range("J53").Value = 9500
$ws.Range("L53").Value = 9500
$ws.Range("N53").Value = -10536

$ws.Range("H55").Value = 409.26086
$ws.Range("I55").Value = 137
$ws.Range("J55").Value = 528.375
$ws.Range("K55").Value = 137
$ws.Range("L55").Value = 528.375
$ws.Range("M55").Value = 36
$ws.Range("N55").Value = -874.375

$ws.Range("H132").Value = 9265380
$ws.Range("I132").Value = 18521380
$ws.Range("K132").Value = 55564140
$ws.Range("M132").Value = -55561610


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H53").Value = 11500
$ws.Range("J53").Value = 11500
$ws.Range("L53").Value = 11500
$ws.Range("N53").Value = -12714

$ws.Range("H122").Value = 204956.2
$ws.Range("I122").Value = 403171
$ws.Range("J122").Value = 6741.4
$ws.Range("K122").Value = 1209513
$ws.Range("L122").Value = 20224.2
$ws.Range("M122").Value = -1207063
$ws.Range("N122").Value = -25124.2

